# Refitting NCDEs to individual patients (for manuscript figure)
# Adds a new "Label" column (H) indicating Control (0) / MDD (1),
# and updates refitted prediction/error/cross-entropy values in columns D, E, F.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add new "Label" header in H1, matching the style of the other headers (G1) ---
$ws.Range("G1").Copy() | Out-Null
$ws.Range("H1").PasteSpecial(-4122) | Out-Null
$ws.Range("H1").Value = "Label"

# --- Updated D/E/F values from the refitted model ---
$ws.Range("D2").Value = 0.893930603357842
$ws.Range("E2").Value = 0.893930603357842

$ws.Range("D4").Value = 0.6937474565470521
$ws.Range("E4").Value = 0.6937474565470521

$ws.Range("D5").Value = 0.6221354905504317
$ws.Range("E5").Value = 0.6221354905504317

$ws.Range("D6").Value = 0.7507661520929523
$ws.Range("E6").Value = 0.7507661520929523

$ws.Range("D7").Value = 0.4098861189543966
$ws.Range("E7").Value = 0.5901138810456035

$ws.Range("D8").Value = 0.6712246356143166
$ws.Range("E8").Value = 0.3287753643856834

$ws.Range("D9").Value = 0.8616438157860511
$ws.Range("E9").Value = 0.1383561842139489

$ws.Range("D10").Value = 0.8025718920574849
$ws.Range("E10").Value = 0.1974281079425151

$ws.Range("F11").Value = 0.8918743133544922

$ws.Range("F21").Value = 0.8681358098983765

# --- New "Label" column values (0 = Control, 1 = MDD) for rows 2-21 ---
$labels = @{
    2 = 0; 3 = 0; 4 = 0; 5 = 0; 6 = 0; 7 = 1; 8 = 1; 9 = 1; 10 = 1; 11 = 1;
    12 = 0; 13 = 0; 14 = 0; 15 = 0; 16 = 0; 17 = 1; 18 = 1; 19 = 1; 20 = 1; 21 = 1;
}

foreach ($row in $labels.Keys) {
    $ws.Range("H$row").Value = $labels[$row]
}
